# Add this week's workout rows (77-82) to the Kilimanjaro weekly scoreboard.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows: Participant, Date (serial), Workout Type, Total Duration,
# Total Distance, Total Elevation, Zone1-5, Workout Level, Week
$newRows = @(
    @{ A = "Eric";     B = 45463; C = "Workout"; D = 80; E = 0;    F = 0;   G = 10; H = 50; I = 15; J = 6;  K = 0; L = "Brave Leopard";  M = 2 },
    @{ A = "Phil";     B = 45463; C = "Run";      D = 74; E = 6.04; F = 928; G = 0;  H = 21; I = 33; J = 11; K = 0; L = "Agile Antelope"; M = 2 },
    @{ A = "Steven";   B = 45463; C = "Walk";     D = 23; E = 1.02; F = 23;  G = 23; H = 0;  I = 0;  J = 0;  K = 0; L = "Agile Antelope"; M = 2 },
    @{ A = "Jeremiah"; B = 45464; C = "Run";      D = 11; E = 1.13; F = 125; G = 0;  H = 4;  I = 2;  J = 0;  K = 0; L = "Agile Antelope"; M = 2 },
    @{ A = "Matt";     B = 45464; C = "Run";      D = 56; E = 6.5;  F = 377; G = 0;  H = 21; I = 29; J = 3;  K = 0; L = "Agile Antelope"; M = 2 },
    @{ A = "Matt";     B = 45464; C = "Walk";     D = 2;  E = 0.15; F = 0;   G = 2;  H = 0;  I = 0;  J = 0;  K = 0; L = "Agile Antelope"; M = 2 }
)

$startRow = 77
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    # Copy the date style (numFmtId 14, "m/d/yyyy") from the row above onto the
    # new date cell so the new rows reuse the existing cell style instead of
    # minting a new one.
    $ws.Range("B" + ($r - 1)).Copy($ws.Range("B" + $r))

    $ws.Range("A$r").Value = $row.A
    $ws.Range("B$r").Value = $row.B
    $ws.Range("C$r").Value = $row.C
    $ws.Range("D$r").Value = $row.D
    $ws.Range("E$r").Value = $row.E
    $ws.Range("F$r").Value = $row.F
    $ws.Range("G$r").Value = $row.G
    $ws.Range("H$r").Value = $row.H
    $ws.Range("I$r").Value = $row.I
    $ws.Range("J$r").Value = $row.J
    $ws.Range("K$r").Value = $row.K
    $ws.Range("L$r").Value = $row.L
    $ws.Range("M$r").Value = $row.M
}

# Move the frozen-pane view / selection the way the author left it after
# scrolling down to the newly entered rows (pane stays frozen at row 1;
# just scroll the bottom pane so row 55 is the first visible row, and
# leave the final selection on the first empty row below the table).
$excel.ActiveWindow.ScrollRow = 55
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("A83").Select()
